$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lab 6 entry gains a descriptive title
$ws.Range("C19").Value = "Lab 6:  Q & A on HW 5 & Missing Data"

# Shorten lecture 14 / 15 labels to draft placeholders (old text kept in column H)
$ws.Range("C26").Value = "Lec 14:"
$ws.Range("C27").Value = "Lec 15: "

# Lecture 12 draft content replaces the old "multiple testing" placeholder
$ws.Range("C20").Value = "Lec 12: Choice of Priors in  Regression"
$ws.Range("D20").Value = "12-reading.html"
$ws.Range("E20").Value = "12-priors-regressions"

# Move the active selection to E20 to match the author's cursor position
$ws.Range("E20").Select()
